$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the snippet id reference for the Field/FieldCollection rows
# (word-manage-fields -> word-document-manage-fields)
$ws.Range("E37").Value = "word-document-manage-fields"
$ws.Range("E38").Value = "word-document-manage-fields"
$ws.Range("E39").Value = "word-document-manage-fields"
$ws.Range("E40").Value = "word-document-manage-fields"
$ws.Range("E41").Value = "word-document-manage-fields"

# Move the active selection to E41
$ws.Range("E41").Select()
